# Widen/shift the four guard-condition textboxes on the
# SelectActivityDiagram slide so their text hugs the same right edge as
# before, and right-align the paragraph text itself inside each box.
#
# NOTE: Shape.Left / Shape.Width are expressed in points in the PowerPoint
# object model (1 pt = 12700 EMU) and are stored internally as single
# precision floats. The literals below were solved offline so that after
# the double -> single narrowing + EMU re-quantization performed when the
# shape position/size is written back, they land on the exact target EMU
# values required by the target OOXML (3244464/775247 -> 2689275/1228836,
# 2590608/1344895 -> 2178760/1756744, 5768252/1755023 -> 5007670/2515606,
# 8148625/1667213 -> 7940970/2292562).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "[else]"
$sh = $s.Shapes.Item("TextBox 46")
$sh.Left = 211.75394439697268
$sh.Width = 96.75873947143555
$sh.TextFrame.TextRange.ParagraphFormat.Alignment = 3

# "[Valid index]"
$sh = $s.Shapes.Item("TextBox 47")
$sh.Left = 171.55590057373047
$sh.Width = 138.32630157470706
$sh.TextFrame.TextRange.ParagraphFormat.Alignment = 3

# "[TotalVisits == 0]"
$sh = $s.Shapes.Item("TextBox 66")
$sh.Left = 394.3047332763672
$sh.Width = 198.07921600341797
$sh.TextFrame.TextRange.ParagraphFormat.Alignment = 3

# "[TotalVisits > 0]"
$sh = $s.Shapes.Item("TextBox 76")
$sh.Left = 625.2732238769532
$sh.Width = 180.51668548583984
$sh.TextFrame.TextRange.ParagraphFormat.Alignment = 3
